# "prep for grafton team"
# Shifts the weekly-route notes in columns E:G (rows 40-45), N/R (rows 39-42),
# and V (rows 38-60), plus the crew-roster blocks in U:W (rows 50-52, 61-63)
# down/along by one slot to make room for the new Grafton team entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column V notes, rows 38-49: shift every value down one row, clear the top ---
$ws.Range("V38").Value = ""
$ws.Range("V39").Value = "5:00 AM MEET OFFICE"
$ws.Range("V40").Value = "6:00 AM START"
$ws.Range("V41").Value = "DC5-FINANCIAL"
$ws.Range("V42").Value = "ROTE OIL #14 TREVOR (CITGO)"
$ws.Range("V43").Value = "12617 ANTIOCH RD"
$ws.Range("V44").Value = "https://goo.gl/maps/xA9YMzPGc6Vhxi5r8"
$ws.Range("V45").Value = "TO FOLLOW"
$ws.Range("V46").Value = "DC5-FINANCIAL"
$ws.Range("V47").Value = "ROTE OIL #13 TREVOR (BP)"
$ws.Range("V48").Value = "12511 ANTIOCH RD"
$ws.Range("V49").Value = "https://goo.gl/maps/iu3xzNwgJ32esP5RA"

# --- Column R, rows 39-41: shift down one row, clear the top ---
$ws.Range("R39").Value = ""
$ws.Range("R40").Value = "5:00 AM OFFICE LEAVE TIME"
$ws.Range("R41").Value = "6:00 AM START"

# --- Column N, rows 40-41: shift down one row, clear the top ---
$ws.Range("N40").Value = ""
$ws.Range("N41").Value = "5:30 AM START "

# --- Columns E:G, rows 40-45: shift every row's content UP one row, clear the bottom ---
$ws.Range("E40").Value = "2)"
$ws.Range("F40").Value = "Carlie"
$ws.Range("G40").Value = "@ Store`n(w/ Trevor)"

$ws.Range("E41").Value = "3)"
$ws.Range("F41").Value = "Jerry D"
$ws.Range("G41").Value = ""

$ws.Range("E42").Value = "4)"
$ws.Range("F42").Value = "Lashaun"

$ws.Range("E43").Value = "5)"
$ws.Range("F43").Value = "Sue"

$ws.Range("E44").Value = "6)"
$ws.Range("F44").Value = "Trevor"
$ws.Range("G44").Value = "@ Store`n(w/ Carlie)"

$ws.Range("E45").Value = ""
$ws.Range("F45").Value = ""
$ws.Range("G45").Value = ""

# --- U:W roster block, rows 50-52: shift down one row, clear the top ---
$ws.Range("U50").Value = ""
$ws.Range("V50").Value = ""
$ws.Range("W50").Value = ""

$ws.Range("U51").Value = "1)"
$ws.Range("V51").Value = "Katherine"
$ws.Range("W51").Value = "Equip"

$ws.Range("U52").Value = "2)"
$ws.Range("V52").Value = "Ian"
$ws.Range("W52").Value = "Driver,`nOptima"

# --- Column V, rows 54-60: shift down one row, clear the top ---
$ws.Range("V54").Value = ""
$ws.Range("V55").Value = "5:30 AM MEET BAYSHORE"
$ws.Range("V56").Value = "7:00 AM START"
$ws.Range("V57").Value = "DC5-ITEM LEVEL"
$ws.Range("V58").Value = "AURORA OUTPATIENT RX #1155 DEPERE"
$ws.Range("V59").Value = "1881 CHICAGO ST"
$ws.Range("V60").Value = "https://goo.gl/maps/9vepkQciLjQDFTRi8"

# --- U:W roster block, rows 61-63: shift down one row, clear the top ---
$ws.Range("U61").Value = ""
$ws.Range("V61").Value = ""
$ws.Range("W61").Value = ""

$ws.Range("U62").Value = "1)"
$ws.Range("V62").Value = "DJ"
$ws.Range("W62").Value = "Driver,`nAltima, Equip"

$ws.Range("U63").Value = "2)"
$ws.Range("V63").Value = "Casey"
$ws.Range("W63").Value = "3rd Aurora"
